{"js": "// The diff:\n//  1) Removes paragraph-level justification (w:jc val=\"both\") from the five\n//     \"SOLID principle\" paragraphs (SRP, OCP, LSP, ISP, DIP).\n//  2) Drops the stray leading space in the run that starts\n//     \" Setiap class memiliki satu tanggung jawab utama. Contohnya, class \"\n//     (inside the SRP paragraph).\n//  3) Fixes a typo \"ooking ID\" -> \"Booking ID\" in the sorting/searching\n//     paragraph (this paragraph's justification is left untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove justified alignment from the five SOLID-principle paragraphs.\nconst principleTitles = [\n  \"Single Responsibility Principle (SRP)\",\n  \"Open/Closed Principle (OCP)\",\n  \"Liskov Substitution Principle (LSP)\",\n  \"Interface Segregation Principle (ISP)\",\n  \"Dependency Inversion Principle (DIP)\",\n];\n\nlet srpParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text;\n  for (const title of principleTitles) {\n    if (text.indexOf(title) === 0) {\n      paragraph.alignment = Word.Alignment.left;\n      if (title === principleTitles[0]) {\n        srpParagraph = paragraph;\n      }\n      break;\n    }\n  }\n}\nawait context.sync();\n\n// 2) Remove the leading space before \"Setiap class memiliki ...\".\nconst leadingSpacePhrase =\n  \" Setiap class memiliki satu tanggung jawab utama. Contohnya, class \";\nconst fixedPhrase =\n  \"Setiap class memiliki satu tanggung jawab utama. Contohnya, class \";\n\nconst srpScope = srpParagraph ? srpParagraph.getRange() : body.getRange();\nconst spaceResults = srpScope.search(leadingSpacePhrase, { matchCase: true });\nspaceResults.load(\"items\");\nawait context.sync();\n\nif (spaceResults.items.length > 0) {\n  spaceResults.items[0].insertText(fixedPhrase, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Fix the \"ooking ID\" typo -> \"Booking ID\".\nconst typoResults = body.search(\"ooking ID\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\n\nif (typoResults.items.length > 0) {\n  typoResults.items[0].insertText(\"Booking ID\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The diff:\n#  1) Removes paragraph-level justification (w:jc val=\"both\") from the five\n#     \"SOLID principle\" paragraphs (SRP, OCP, LSP, ISP, DIP).\n#  2) Drops the stray leading space in the run that starts\n#     \" Setiap class memiliki satu tanggung jawab utama. Contohnya, class \"\n#     (inside the SRP paragraph).\n#  3) Fixes a typo \"ooking ID\" -> \"Booking ID\" in the sorting/searching\n#     paragraph (this paragraph's justification is left untouched).\n\n$d = $word.ActiveDocument\n\n# 1) Remove justified alignment from the five SOLID-principle paragraphs.\n$principleTitles = @(\n    \"Single Responsibility Principle (SRP)\",\n    \"Open/Closed Principle (OCP)\",\n    \"Liskov Substitution Principle (LSP)\",\n    \"Interface Segregation Principle (ISP)\",\n    \"Dependency Inversion Principle (DIP)\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    foreach ($title in $principleTitles) {\n        if ($t.StartsWith($title)) {\n            $p.Format.Alignment = 0   # wdAlignParagraphLeft\n            break\n        }\n    }\n}\n\n# 2) Remove the leading space before \"Setiap class memiliki ...\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" Setiap class memiliki satu tanggung jawab utama. Contohnya, class \"\nif ($find.Execute()) {\n    $find.Parent.Text = \"Setiap class memiliki satu tanggung jawab utama. Contohnya, class \"\n}\n\n# 3) Fix the \"ooking ID\" typo -> \"Booking ID\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"ooking ID\"\nif ($find2.Execute()) {\n    $find2.Parent.Text = \"Booking ID\"\n}\n"}
